# Updates cryptocurrency Price (D) and Volume(1h) (E) columns with refreshed
# market data, matching the GitHub Actions scheduled refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (no border/shading/number-format) used to restore cell
# styling after the quote-prefix trick below forces Price values to stay text.
$plainStyle = $ws.Range("A1").Style

$ws.Range("D2").Value = "'24.900.12"
$ws.Range("D2").Style = $plainStyle
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "'1.700.13"
$ws.Range("D3").Style = $plainStyle
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = $plainStyle
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("D5").Value = "'315.22"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").Value = "'0.4025"
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = "  +2.32%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'1.004"
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").Value = "'53.70"
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("E11").Value = "  -3.46%  "
$ws.Range("D12").Value = "'0.08816"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "'25.77"
$ws.Range("D13").Style = $plainStyle
$ws.Range("E13").Value = "  +5.16%  "
$ws.Range("D14").Value = "'7.485"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").Value = "'8.034"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").Value = "'0.00001346"
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("D17").Value = "'1.817.37"
$ws.Range("D17").Style = $plainStyle
$ws.Range("E17").Value = "  +7.46%  "
$ws.Range("D18").Value = "'96.68"
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("D19").Value = "'0.07179"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("D20").Value = "'20.96"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "  +5.74%  "
$ws.Range("D21").Value = "'7.235"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  -2.21%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("D23").Value = "'14.59"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("D24").Value = "'24.896.50"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("D25").Value = "'2.339"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("D26").Value = "'2.889"
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = "  -5.45%  "
$ws.Range("D27").Value = "'6.596"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  +26.52%  "
$ws.Range("D28").Value = "'23.06"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("D29").Value = "'163.58"
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("D30").Value = "'143.54"
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = "  +4.14%  "
$ws.Range("D31").Value = "'8.155"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = "  -3.85%  "
$ws.Range("D32").Value = "'1.999.01"
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = "  +6.34%  "
$ws.Range("D33").Value = "'2.279"
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = "  +14.12%  "
$ws.Range("D34").Value = "'0.08758"
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = "  -1.24%  "
$ws.Range("D35").Value = "'7.405"
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("D36").Value = "'0.03184"
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = "  +9.15%  "
$ws.Range("D37").Value = "'1.032"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("D39").Value = "'0.8510"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  +8.29%  "
$ws.Range("D40").Value = "'10.91"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("D41").Value = "'0.09424"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  +2.97%  "
$ws.Range("D42").Value = "'14.03"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = "  -2.30%  "
$ws.Range("D43").Value = "'1.470"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").Value = "'17.77"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  +6.70%  "
$ws.Range("D45").Value = "'2.715"
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = "  +5.33%  "
$ws.Range("D46").Value = "'0.7467"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  +3.36%  "
$ws.Range("D47").Value = "'4.246"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("D48").Value = "'1.408"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  +5.72%  "
$ws.Range("D49").Value = "'1.005"
$ws.Range("D49").Style = $plainStyle
$ws.Range("D50").Value = "'141.37"
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("D51").Value = "'0.08352"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  +4.52%  "
